# 10th - MB for single stock and added new group
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data columns (Jun_27, Jun_26) right after the firm-name column A.
# This shifts the existing Jun_10 data from column B to column D.
$ws.Range("B1:C1").EntireColumn.Insert()

# New column headers
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"

# Fill in "UN" (unchanged) placeholder for every existing firm row in the two new columns
$ws.Range("B2:C27").Value = "UN"

# This week's (Jun_27) new rating: Morgan Stanley (row 7)
$ws.Range("B7").Value = "6/27/2018,Reiterates,Equal Weight -> Overweight,$37.00"

# Prior week's (Jun_26) new rating: JPMorgan Chase & Co. (row 17)
$ws.Range("C17").Value = "6/19/2018,Initiates,Overweight -> Neutral,$36.00"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 21.1640625
$ws.Columns.Item(2).ColumnWidth = 39.83203125
$ws.Columns.Item(3).ColumnWidth = 39.83203125
$ws.Columns.Item(4).ColumnWidth = 61.1640625

# New analyst groups added at the bottom of the table
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"

$ws.Range("B11").Select()
